# Insert a new data row immediately before the current row 59, shifting the
# existing rows 59:162 down to 60:163 (values, formatting and styles move
# with them, exactly like Excel's own Rows.Insert()).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(59).Insert()

# Populate the freshly inserted row 59 with the new record.
$ws.Range("A59").Value = 7
$ws.Range("B59").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C59").Value = "Ñuble"
$ws.Range("D59").Value = 44757
$ws.Range("E59").Value = 16
$ws.Range("F59").Value = 100112045
$ws.Range("G59").Value = "Zapallo"
$ws.Range("H59").Value = "Camote"
$ws.Range("I59").Value = "1a (guarda)"
$ws.Range("J59").Value = 120
$ws.Range("K59").Value = 550
$ws.Range("L59").Value = 600
$ws.Range("M59").Value = 575
$ws.Range("N59").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O59").Value = "Región del Maule"
$ws.Range("P59").Value = 575
$ws.Range("Q59").Value = 1
$ws.Range("R59").Value = "Hortaliza"
